$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.896.97'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '2.354.23'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +5.90%  '
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '2.774.17'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = '57.832.75'
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '2.360.26'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '330.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("E21").Value = '  -2.54%  '
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '62.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("E27").Value = '  -3.23%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").Value = '0.0₃0735'
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  -2.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '142.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.378'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '288.94'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0222'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.381'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.64%  '
